{"js": "const replacements = [\n  [\"672\u00f78=\", \"602\u00f75=\"],\n  [\"885\u00f73=\", \"967\u00f74=\"],\n  [\"759\u00f76=\", \"902\u00f79=\"],\n  [\"952\u00f77=\", \"927\u00f77=\"],\n  [\"961\u00f72=\", \"677\u00f75=\"],\n  [\"449\u00f74=\", \"291\u00f79=\"],\n  [\"561\u00f75=\", \"966\u00f74=\"],\n  [\"893\u00f79=\", \"308\u00f77=\"],\n  [\"937\u00f75=\", \"790\u00f76=\"],\n  [\"661\u00f73=\", \"437\u00f77=\"],\n  [\"676\u00f73=\", \"401\u00f74=\"],\n  [\"660\u00f75=\", \"818\u00f79=\"],\n  [\"436\u00f78=\", \"346\u00f78=\"],\n  [\"931\u00f74=\", \"170\u00f74=\"],\n  [\"984\u00f77=\", \"384\u00f78=\"],\n  [\"288\u00f76=\", \"409\u00f76=\"],\n  [\"873\u00f78=\", \"556\u00f79=\"],\n  [\"965\u00f72=\", \"986\u00f73=\"],\n  [\"101\u00f76=\", \"735\u00f79=\"],\n  [\"334\u00f76=\", \"707\u00f73=\"],\n  [\"457\u00f73=\", \"576\u00f75=\"],\n  [\"245\u00f75=\", \"879\u00f76=\"],\n  [\"353\u00f74=\", \"204\u00f75=\"],\n  [\"793\u00f73=\", \"338\u00f78=\"],\n  [\"452\u00f77=\", \"897\u00f74=\"],\n];\n\nconst body = context.document.body;\n\n// Issue all searches first, then load them together.\nconst searchResults = replacements.map(([from]) =>\n  body.search(from, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((r) => r.load(\"items\"));\nawait context.sync();\n\nlet totalFound = 0;\nsearchResults.forEach((results, idx) => {\n  const [, to] = replacements[idx];\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, Word.InsertLocation.replace);\n    totalFound++;\n  }\n});\nawait context.sync();\n\nreturn \"replaced: \" + totalFound;", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"672\u00f78=\", \"602\u00f75=\"),\n    @(\"885\u00f73=\", \"967\u00f74=\"),\n    @(\"759\u00f76=\", \"902\u00f79=\"),\n    @(\"952\u00f77=\", \"927\u00f77=\"),\n    @(\"961\u00f72=\", \"677\u00f75=\"),\n    @(\"449\u00f74=\", \"291\u00f79=\"),\n    @(\"561\u00f75=\", \"966\u00f74=\"),\n    @(\"893\u00f79=\", \"308\u00f77=\"),\n    @(\"937\u00f75=\", \"790\u00f76=\"),\n    @(\"661\u00f73=\", \"437\u00f77=\"),\n    @(\"676\u00f73=\", \"401\u00f74=\"),\n    @(\"660\u00f75=\", \"818\u00f79=\"),\n    @(\"436\u00f78=\", \"346\u00f78=\"),\n    @(\"931\u00f74=\", \"170\u00f74=\"),\n    @(\"984\u00f77=\", \"384\u00f78=\"),\n    @(\"288\u00f76=\", \"409\u00f76=\"),\n    @(\"873\u00f78=\", \"556\u00f79=\"),\n    @(\"965\u00f72=\", \"986\u00f73=\"),\n    @(\"101\u00f76=\", \"735\u00f79=\"),\n    @(\"334\u00f76=\", \"707\u00f73=\"),\n    @(\"457\u00f73=\", \"576\u00f75=\"),\n    @(\"245\u00f75=\", \"879\u00f76=\"),\n    @(\"353\u00f74=\", \"204\u00f75=\"),\n    @(\"793\u00f73=\", \"338\u00f78=\"),\n    @(\"452\u00f77=\", \"897\u00f74=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n\nWrite-Output \"done\""}
